$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update stats for 2025-09 (row 22)
$ws.Range("B22").Value = 6295
$ws.Range("D22").Value = 5848636
$ws.Range("E22").Value = 929.0922954725972
$ws.Range("F22").Value = 8.366328111551047
$ws.Range("H22").Value = 27.18960672935569
